# Commit: add project_type in Project and calculate score to be loaded in
# MutationSummary as well as variant_caller_presence
#
# Net effect on the workbook: insert a new "project_type" column into the
# "Project" sheet, between the existing "start_date" and "description"
# columns (i.e. the new column becomes column H, and the old "description"
# column shifts from H to I). Also the "Project" sheet becomes the active
# tab/selection (it was "GuideMismatches" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project")

# Insert a new column at H, shifting "description" (and everything after
# it) one column to the right.
$ws.Columns.Item(8).Insert() | Out-Null

# New header for the inserted column.
$ws.Range("H1").Value = "project_type"

# Match the template's column width for the new column (closest value the
# COM width model can represent).
$ws.Columns.Item(8).ColumnWidth = 10.3

# Make "Project" the active sheet/tab, with H2 selected (the cell under
# the new header) as in the edited template.
$ws.Activate() | Out-Null
$ws.Range("H2").Select() | Out-Null
